# Update sector names in column D to match renamed sector taxonomy:
#   "Health Care"            -> "Healthcare"
#   "Consumer Staples"       -> "Consumer Defensive"
#   "Consumer Discretionary" -> "Consumer Cyclical"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "Health Care"            = "Healthcare"
    "Consumer Staples"       = "Consumer Defensive"
    "Consumer Discretionary" = "Consumer Cyclical"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
